# "Add files via upload" — re-label the three monthly header bands
# (Jan..Dec) on row 3 of the condensed raw-data sheet so each block
# carries its year suffix: Jan_2019..Dec_2019 (2019 block, cols B:M),
# Jan_2020..Dec_2020 (2020 block, cols O:Z) and Jan_2021..Nov_2021
# (2021 block, cols AB:AL, only 11 months of data present).
# The underlying numeric data (rows 4-22) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$months = @("Jan", "Feb", "Mar", "Apr", "May", "Jun", "Jul", "Aug", "Sep", "Oct", "Nov", "Dec")

# 2019 block: row 3, columns B..M
$cols2019 = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M")
for ($i = 0; $i -lt $cols2019.Length; $i++) {
    $ws.Range($cols2019[$i] + "3").Value = $months[$i] + "_2019"
}

# 2020 block: row 3, columns O..Z
$cols2020 = @("O", "P", "Q", "R", "S", "T", "U", "V", "W", "X", "Y", "Z")
for ($i = 0; $i -lt $cols2020.Length; $i++) {
    $ws.Range($cols2020[$i] + "3").Value = $months[$i] + "_2020"
}

# 2021 block: row 3, columns AB..AL (Jan through Nov only)
$cols2021 = @("AB", "AC", "AD", "AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL")
for ($i = 0; $i -lt $cols2021.Length; $i++) {
    $ws.Range($cols2021[$i] + "3").Value = $months[$i] + "_2021"
}

# Reflect the author's final selection / scroll position in the saved view.
$ws.Range("AM3").Select()
